# Update the "Comment" (column E) text values used throughout the change log.
# The underlying wording for each category changed; re-point every cell that
# used the old phrase to the corresponding new phrase.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldToNew = @{
    "Network added to Base Service in Mar 2020"          = "Network Added to Base Service in Mar 2020"
    "Network added to Add-On Service in Mar 2020"        = "Network Added to Add-On Package"
    "Network removed from Add-On Service in Mar 2020"    = "Network Removed from Add-On Package"
    "Network moved to new Add-On Service in Mar 2020"    = "Network Moved from One Add-On Package to Another Add-On Package"
    "New Network Added to Database on Mar 2020"          = "New Network Added to Database"
    "Existing Network Removed from Database in Mar 2020" = "Network Removed from Database"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -ne $null -and $oldToNew.ContainsKey($val)) {
        $cell.Value = $oldToNew[$val]
    }
}

$ws.Range("E45:E46").Select()
